$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '59.836.35'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'" + '  -3.35%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = "'" + '3.277.44'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "'" + '  -4.11%  '
$ws.Range('E3').ClearFormats()
$ws.Range('D4').Value = "'" + '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = "'" + '  +0.03%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = "'" + '554.49'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "'" + '  -4.14%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = "'" + '139.89'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "'" + '  -8.70%  '
$ws.Range('E6').ClearFormats()
$ws.Range('E7').Value = "'" + '  -0.05%  '
$ws.Range('E7').ClearFormats()
$ws.Range('D8').Value = "'" + '3.275.05'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = "'" + '  -4.15%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').Value = "'" + '0.464'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "'" + '  -4.10%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').Value = "'" + '7.82'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "'" + '  -3.01%  '
$ws.Range('E10').ClearFormats()
$ws.Range('E11').Value = "'" + '  -5.53%  '
$ws.Range('E11').ClearFormats()
$ws.Range('D12').Value = "'" + '0.403'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "'" + '  -3.64%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').Value = "'" + '3.833.44'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "'" + '  -4.22%  '
$ws.Range('E13').ClearFormats()
$ws.Range('E14').Value = "'" + '  -0.22%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').Value = "'" + '26.56'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = "'" + '  -6.95%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').Value = "'" + '3.263.32'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = "'" + '  -4.13%  '
$ws.Range('E16').ClearFormats()
$ws.Range('E17').Value = "'" + '  -5.16%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').Value = "'" + '59.884.78'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'" + '  -3.36%  '
$ws.Range('E18').ClearFormats()
$ws.Range('D19').Value = "'" + '6.05'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "'" + '  -7.24%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').Value = "'" + '13.65'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = "'" + '  -6.00%  '
$ws.Range('E20').ClearFormats()
$ws.Range('D21').Value = "'" + '8.47'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'" + '  -5.50%  '
$ws.Range('E21').ClearFormats()
$ws.Range('D22').Value = "'" + '371.68'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'" + '  -2.72%  '
$ws.Range('E22').ClearFormats()
$ws.Range('B23').Value = "'" + 'Dai'
$ws.Range('B23').ClearFormats()
$ws.Range('C23').Value = "'" + 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('C23').ClearFormats()
$ws.Range('D23').Value = "'" + '1.00'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = "'" + '  +0.13%  '
$ws.Range('E23').ClearFormats()
$ws.Range('B24').Value = "'" + 'Litecoin'
$ws.Range('B24').ClearFormats()
$ws.Range('C24').Value = "'" + 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('C24').ClearFormats()
$ws.Range('D24').Value = "'" + '72.46'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "'" + '  -3.63%  '
$ws.Range('E24').ClearFormats()
$ws.Range('D25').Value = "'" + '0.528'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "'" + '  -7.59%  '
$ws.Range('E25').ClearFormats()
$ws.Range('D26').Value = "'" + '3.401.63'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "'" + '  -4.42%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').Value = "'" + '0.0000102'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "'" + '  -9.29%  '
$ws.Range('E27').ClearFormats()
$ws.Range('D28').Value = "'" + '0.171'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "'" + '  -4.18%  '
$ws.Range('E28').ClearFormats()
$ws.Range('D29').Value = "'" + '1.00'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = "'" + '  +0.52%  '
$ws.Range('E29').ClearFormats()
$ws.Range('E30').Value = "'" + '  -8.87%  '
$ws.Range('E30').ClearFormats()
$ws.Range('E31').Value = "'" + '  -0.02%  '
$ws.Range('E31').ClearFormats()
$ws.Range('D32').Value = "'" + '2.01'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = "'" + '  -5.35%  '
$ws.Range('E32').ClearFormats()
$ws.Range('D33').Value = "'" + '7.42'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "'" + '  -5.65%  '
$ws.Range('E33').ClearFormats()
$ws.Range('D34').Value = "'" + '22.42'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = "'" + '  -3.31%  '
$ws.Range('E34').ClearFormats()
$ws.Range('D35').Value = "'" + '1.23'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = "'" + '  -7.91%  '
$ws.Range('E35').ClearFormats()
$ws.Range('D36').Value = "'" + '165.98'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = "'" + '  -1.58%  '
$ws.Range('E36').ClearFormats()
$ws.Range('D37').Value = "'" + '5.03'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "'" + '  -8.35%  '
$ws.Range('E37').ClearFormats()
$ws.Range('E38').Value = "'" + '  -5.51%  '
$ws.Range('E38').ClearFormats()
$ws.Range('D39').Value = "'" + '6.60'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "'" + '  -5.28%  '
$ws.Range('E39').ClearFormats()
$ws.Range('D40').Value = "'" + '3.304.43'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'" + '  -4.30%  '
$ws.Range('E40').ClearFormats()
$ws.Range('B41').Value = "'" + 'Hedera'
$ws.Range('B41').ClearFormats()
$ws.Range('C41').Value = "'" + 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C41').ClearFormats()
$ws.Range('D41').Value = "'" + '0.0721'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "'" + '  -8.30%  '
$ws.Range('E41').ClearFormats()
$ws.Range('B42').Value = "'" + 'EnergySwap'
$ws.Range('B42').ClearFormats()
$ws.Range('C42').Value = "'" + 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C42').ClearFormats()
$ws.Range('D42').Value = "'" + '25.60'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = "'" + '  -17.46%  '
$ws.Range('E42').ClearFormats()
$ws.Range('D43').Value = "'" + '41.39'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = "'" + '  -2.90%  '
$ws.Range('E43').ClearFormats()
$ws.Range('D44').Value = "'" + '0.743'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'" + '  -4.85%  '
$ws.Range('E44').ClearFormats()
$ws.Range('B45').Value = "'" + 'ONDO'
$ws.Range('B45').ClearFormats()
$ws.Range('C45').Value = "'" + 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('C45').ClearFormats()
$ws.Range('D45').Value = "'" + '1.12'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'" + '  -3.81%  '
$ws.Range('E45').ClearFormats()
$ws.Range('B46').Value = "'" + 'Filecoin'
$ws.Range('B46').ClearFormats()
$ws.Range('C46').Value = "'" + 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C46').ClearFormats()
$ws.Range('D46').Value = "'" + '4.09'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = "'" + '  -7.60%  '
$ws.Range('E46').ClearFormats()
$ws.Range('D47').Value = "'" + '1.56'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "'" + '  -7.17%  '
$ws.Range('E47').ClearFormats()
$ws.Range('E48').Value = "'" + '  -0.01%  '
$ws.Range('E48').ClearFormats()
$ws.Range('D49').Value = "'" + '2.313.42'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "'" + '  -9.24%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').Value = "'" + '6.31'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "'" + '  -7.72%  '
$ws.Range('E50').ClearFormats()
$ws.Range('D51').Value = "'" + '21.20'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = "'" + '  -5.80%  '
$ws.Range('E51').ClearFormats()
